# "Generate Report for Archive": refresh the localization-status report so
# the in-flight item (currently "Ready for handoff") reflects that it has
# moved into translation, then re-pack the Status column widths the way the
# report generator does whenever the status text changes length.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-locale status columns (E = zh-cn, F = de-de).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale detail sheets: Status column (C).
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Column widths on the Status columns are tightened to fit the new,
# shorter status text.
$newWidth = 13.4101845877511 - 0.8333333333333333

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
